# Natmi following Dr Hou advice: refresh the Wnt5a-Ror2 LR-pair results
# (YoungD7) so FAPs/sCs are both sending AND target clusters, with ECs
# added as an additional target cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One variable per data row (columns A..T); avoids nested array literals.
$dataRow1 = @("FAPs", "Wnt5a", "Ror2", "ECs", 3, 1, 10.43424333333333, 31.30273, 0.9711091978791583, 0.9711091978791584, 1, 0.3333333333333333, 0.04009133333333333, 0.120274, 0.004449173663854228, 0.004449173663854228, 0.4183227275577778, 3.76490454802, 0.004320633467930554, 0.004320633467930556)
$dataRow2 = @("FAPs", "Wnt5a", "Ror2", "FAPs", 3, 1, 10.43424333333333, 31.30273, 0.9711091978791583, 0.9711091978791584, 3, 1, 8.212806333333333, 24.638419, 0.9114239564145669, 0.911423956414567, 85.6944197315411, 771.2497775838699, 0.885092187241599, 0.8850921872415991)
$dataRow3 = @("FAPs", "Wnt5a", "Ror2", "sCs", 3, 1, 10.43424333333333, 31.30273, 0.9711091978791583, 0.9711091978791584, 3, 1, 0.7580640000000001, 2.274192, 0.08412686992157885, 0.08412686992157886, 7.909824238240001, 71.18841814416, 0.08169637716962873, 0.08169637716962874)
$dataRow4 = @("sCs", "Wnt5a", "Ror2", "ECs", 2, 0.6666666666666666, 0.310422, 0.9312659999999999, 0.02889080212084161, 0.02889080212084161, 1, 0.3333333333333333, 0.04009133333333333, 0.120274, 0.004449173663854228, 0.004449173663854228, 0.012445231876, 0.112007086884, 0.0001285401959236723, 0.0001285401959236724)
$dataRow5 = @("sCs", "Wnt5a", "Ror2", "FAPs", 2, 0.6666666666666666, 0.310422, 0.9312659999999999, 0.02889080212084161, 0.02889080212084161, 3, 1, 8.212806333333333, 24.638419, 0.9114239564145669, 0.911423956414567, 2.549435767606, 22.944921908454, 0.02633176917296782, 0.02633176917296782)
$dataRow6 = @("sCs", "Wnt5a", "Ror2", "sCs", 2, 0.6666666666666666, 0.310422, 0.9312659999999999, 0.02889080212084161, 0.02889080212084161, 3, 1, 0.7580640000000001, 2.274192, 0.08412686992157885, 0.08412686992157886, 0.235319743008, 2.117877687072, 0.002430492751950116, 0.002430492751950117)

$rows = @($dataRow1, $dataRow2, $dataRow3, $dataRow4, $dataRow5, $dataRow6)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $sheetRow = $i + 2
    $row = $rows[$i]
    for ($col = 1; $col -le $row.Count; $col++) {
        $ws.Cells.Item($sheetRow, $col).Value = $row[$col - 1]
    }
}
